# This script migrates the workbook from the newer layout (with a "Texas Notes"
# worksheet and its associated notes/sources) back to the earlier layout used by
# the "3.0 RMI script" data drop: the "Texas Notes" sheet is removed, the
# EoDSDwSP sheet's formulas are repointed at the Calculations sheet (instead of
# the now-removed Texas Notes sheet), and the old hyperlink styling on the
# About sheet is cleared.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsCalc  = $wb.Worksheets.Item("Calculations")
$wsEoD   = $wb.Worksheets.Item("EoDSDwSP")

# --- EoDSDwSP: repoint formulas away from the soon-to-be-deleted "Texas Notes"
# sheet so they resolve against Calculations instead (avoids #REF! errors once
# the sheet is removed below).
$wsEoD.Range("B2").Formula = "=Calculations!B9"
$wsEoD.Range("B4").Formula = "=Calculations!B10"

# --- About: remove the hyperlink and its "Hyperlink" styling from B6, leaving
# plain text behind.
$rngLink = $wsAbout.Range("B6")
$rngLink.Hyperlinks.Delete()
$rngLink.Style = "Normal"

# Drop the now-unused built-in "Hyperlink" cell style definition.
$wb.Styles.Item("Hyperlink").Delete()

# --- Remove the "Texas Notes" worksheet entirely.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Texas Notes").Delete()
$excel.DisplayAlerts = $true

# Re-fetch worksheet references after the deletion above, since sheet
# handles/indices shift once a sheet is removed.
$wsAbout = $wb.Worksheets.Item("About")
$wsCalc  = $wb.Worksheets.Item("Calculations")
$wsEoD   = $wb.Worksheets.Item("EoDSDwSP")

# --- Restore the selections/active cells used in the earlier workbook layout.
$wsCalc.Activate()
$wsCalc.Range("A1").Select()

$wsEoD.Activate()
$wsEoD.Range("B2").Select()

$wsAbout.Activate()
$wsAbout.Range("A12").Select()
